$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block updates (plain text, no numeric coercion risk)
$ws.Range("B2").Value = "12345-ABC"
$ws.Range("E2").Value = "adasdas dasd 66 Q"

# Role specialization update
$ws.Range("B7").Value = "DevOps Engineer"

# Jan 02 row: move 1.0 from Sick Leave (E12) to At Work (C12)
# Force text storage ("1.0") instead of Excel's automatic numeric coercion,
# matching the original workbook's text-based numeric cells.
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value2 = "1.0"
$ws.Range("E12").Value = ""

# Jan 03 row: move 1.0 from Sick Leave (E13) to At Work (C13)
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value2 = "1.0"
$ws.Range("E13").Value = ""

# Totals row updates
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value2 = "20.0"
$ws.Range("E44").Value = "-"

# Signature date update
$ws.Range("B50").Value = "08 - February - 2025"
